# Generate Report for Handback
# Adds a new handback entry (8e6706c6-cdb4-4fb6-b910-12cd42910fde.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the
# shape of the two existing rows (c0f7397b... / 26e4cf4a...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "8e6706c6-cdb4-4fb6-b910-12cd42910fde.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-10-17 14:16:53"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e6706c6cdb44fb6b91012cd42910fde/e2e/8e6706c6-cdb4-4fb6-b910-12cd42910fde.md",
    "",
    "",
    "e2e\8e6706c6-cdb4-4fb6-b910-12cd42910fde.md"
) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = "8e6706c6-cdb4-4fb6-b910-12cd42910fde.1b54bea508832d58774593abbf6dc4c27c72d84f.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-10-17 14:16:30"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("J4").Value = "8e6706c6-cdb4-4fb6-b910-12cd42910fde.1b54bea508832d58774593abbf6dc4c27c72d84f.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-10-17 14:17:35"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e6706c6cdb44fb6b91012cd42910fde/e2e/8e6706c6-cdb4-4fb6-b910-12cd42910fde.md",
    "",
    "",
    "8e6706c6-cdb4-4fb6-b910-12cd42910fde.md"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8e6706c6cdb44fb6b91012cd42910fde/e2e/8e6706c6-cdb4-4fb6-b910-12cd42910fde.md",
    "",
    "",
    "8e6706c6-cdb4-4fb6-b910-12cd42910fde.md"
) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = "8e6706c6-cdb4-4fb6-b910-12cd42910fde.1b54bea508832d58774593abbf6dc4c27c72d84f.de-de.xlf"
$wsDe.Range("H4").Value = "2016-10-17 14:16:53"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("J4").Value = "8e6706c6-cdb4-4fb6-b910-12cd42910fde.1b54bea508832d58774593abbf6dc4c27c72d84f.de-de.xlf"
$wsDe.Range("K4").Value = "2016-10-17 14:18:12"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e6706c6cdb44fb6b91012cd42910fde/e2e/8e6706c6-cdb4-4fb6-b910-12cd42910fde.md",
    "",
    "",
    "8e6706c6-cdb4-4fb6-b910-12cd42910fde.md"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8e6706c6cdb44fb6b91012cd42910fde/e2e/8e6706c6-cdb4-4fb6-b910-12cd42910fde.md",
    "",
    "",
    "8e6706c6-cdb4-4fb6-b910-12cd42910fde.md"
) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))
